$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A1: becomes "Bold Times New Roman" in a bold Times New Roman font ---
# (it keeps using the same style slot the cell already occupied, now redefined)
$ws.Range("A1").Value = "Bold Times New Roman"
$ws.Range("A1").Font.Italic = $false
$ws.Range("A1").Font.Size = 11
$ws.Range("A1").Font.Name = "Times New Roman"
$ws.Range("A1").Font.Bold = $true

# --- B3: new cell carrying the original 24pt italic font forward ---
$ws.Range("B3").Value = "24 pt Italic"
$ws.Range("B3").Font.Italic = $true
$ws.Range("B3").Font.Size = 24
